$newValues = @(
    "69+28=",
    "95-79=",
    "11-3=",
    "25-17=",
    "43-27=",
    "35+57=",
    "82-13=",
    "63-38=",
    "5+69=",
    "37+9=",
    "14+38=",
    "28+66=",
    "38+46=",
    "17+55=",
    "22+69=",
    "93-9=",
    "18+16=",
    "17+8=",
    "87-59=",
    "23-4=",
    "73+9=",
    "66-7=",
    "65+8=",
    "80-4=",
    "82-69=",
    "27+66=",
    "95-58=",
    "15+19=",
    "61-18=",
    "90-13=",
    "45+16=",
    "44+8=",
    "49+13=",
    "62-27=",
    "8+7=",
    "7+28=",
    "61-24=",
    "70-15=",
    "5+37=",
    "22-6=",
    "63-49=",
    "71-7=",
    "58+34=",
    "73-15=",
    "52-35=",
    "9+7=",
    "12-8=",
    "6+38=",
    "47+24=",
    "73-47=",
    "92-87=",
    "55-19=",
    "72-55=",
    "2+29=",
    "73-65=",
    "39+17=",
    "59+2=",
    "4+28=",
    "29+22=",
    "69+23=",
    "15+47=",
    "35-26=",
    "72-46=",
    "74-27=",
    "48+44=",
    "93-88=",
    "94-87=",
    "23+18=",
    "46+15=",
    "76+19=",
    "16+69=",
    "70-11=",
    "86-58=",
    "30-6=",
    "19+22=",
    "16+19=",
    "33-15=",
    "69+6=",
    "8+53=",
    "45-37=",
    "90-2=",
    "18+48=",
    "19+44=",
    "40-11=",
    "46+45=",
    "75-69=",
    "46+48=",
    "55-9=",
    "44+39=",
    "17+46=",
    "45+36=",
    "30-19=",
    "28+56=",
    "63-18=",
    "8+37=",
    "83-18=",
    "84-46=",
    "34-7=",
    "28-9=",
    "75+9="
)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Output ("Updated " + $idx + " cells")